$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns, matching the formatting of column H's header.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for rows 2-14.
$values = @(
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(6, 7),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
